$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): names get a trailing newline, "Point" labels get a colon ---
$ws.Range("A1").Value = "Kim`n"
$ws.Range("B1").Value = "Point:"
$ws.Range("C1").Value = "Emil`n"
$ws.Range("D1").Value = "Point:"
$ws.Range("E1").Value = "Mads`n"
$ws.Range("F1").Value = "Point:"
$ws.Range("G1").Value = "Soren`n"
$ws.Range("H1").Value = "Point:"

# Re-fit row 1 height so the embedded newline doesn't leave a custom row height behind
$ws.Rows.Item(1).AutoFit()

# --- Column A (Kim) ---
$ws.Range("A3").Value = "RB leipzig"
$ws.Range("A4").Value = "Barcelona"
$ws.Range("A5").Value = "Bologna"
$ws.Range("A6").Value = "OB"
$ws.Range("A7").Value = "AGF"

# --- Column C (Emil) ---
$ws.Range("C2").Value = "Leicester"
$ws.Range("C3").Value = "Tottenham"
$ws.Range("C4").Value = "Dortmund"
$ws.Range("C5").Value = "Real sociedad"
$ws.Range("C6").Value = "Atalanta"
$ws.Range("C7").Value = "Brndby IF"

# --- Column E (Mads) ---
$ws.Range("E2").Value = "Arsenal"
$ws.Range("E3").Value = "Frankfurt"
$ws.Range("E4").Value = "Hoffenheim"
$ws.Range("E5").Value = "Valencia"
$ws.Range("E6").Value = "Ac Milan"
$ws.Range("E7").Value = "FC Kbenhavn"

# --- Column G (Soren) ---
$ws.Range("G2").Value = "Manchester utd"
$ws.Range("G3").Value = "Leverkusen"
$ws.Range("G4").Value = "sevilla"
$ws.Range("G5").Value = "juventus"
$ws.Range("G6").Value = "Torino"
$ws.Range("G7").Value = "Fc midtjylland"
